# Updated symbol list on Wed Jan 18 22:20:36 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) readings for
# the coin rows on the sheet. Values are stored as plain text (e.g. "289.13",
# "-4.28%") rather than numbers, so each cell is forced to Text format before
# the write (and reset back to the default "Normal" style right after) to
# stop Excel's COM layer from auto-coercing the literal into a float/percent
# and silently dropping significant trailing zeros (e.g. "0.8990" -> 0.899).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "289.13"
Set-TextValue "E2" "-4.28%"

Set-TextValue "D3" "30.84"
Set-TextValue "E3" "-4.09%"

Set-TextValue "D4" "4.893"
Set-TextValue "E4" "-1.89%"

Set-TextValue "D5" "0.07152"
Set-TextValue "E5" "-9.57%"

Set-TextValue "D6" "1.833"
Set-TextValue "E6" "-12.73%"

Set-TextValue "D7" "7.653"
Set-TextValue "E7" "-2.19%"

Set-TextValue "D8" "3.725"
Set-TextValue "E8" "-2.00%"

Set-TextValue "D9" "0.8990"
Set-TextValue "E9" "-2.94%"

Set-TextValue "D10" "0.1648"
Set-TextValue "E10" "-5.69%"

Set-TextValue "D11" "0.07512"
Set-TextValue "E11" "-5.54%"

Set-TextValue "D12" "0.08122"
Set-TextValue "E12" "-6.06%"

Set-TextValue "D13" "0.02994"
Set-TextValue "E13" "-4.53%"

Set-TextValue "D14" "0.09998"
Set-TextValue "E14" "-0.20%"

Set-TextValue "D15" "0.001507"
Set-TextValue "E15" "-0.18%"

Set-TextValue "D16" "0.005704"
Set-TextValue "E16" "-1.34%"

Set-TextValue "D18" "3.460"
Set-TextValue "E18" "-0.06%"

Set-TextValue "D19" "2.104"
Set-TextValue "E19" "-7.52%"

Set-TextValue "D20" "0.3276"
Set-TextValue "E20" "-0.33%"

Set-TextValue "D21" "0.1299"
Set-TextValue "E21" "0.68%"

Set-TextValue "D22" "4.370"
Set-TextValue "E22" "1.22%"

Set-TextValue "D23" "0.2004"
Set-TextValue "E23" "11.90%"

Set-TextValue "D24" "0.04468"
Set-TextValue "E24" "-2.91%"

Set-TextValue "E25" "-2.08%"

Set-TextValue "D26" "0.004018"
Set-TextValue "E26" "-9.99%"

Set-TextValue "D27" "0.0001253"
Set-TextValue "E27" "0.18%"

Set-TextValue "D39" "0.01640"
Set-TextValue "E39" "-4.81%"

Set-TextValue "D40" "0.04332"
Set-TextValue "E40" "-9.47%"

Set-TextValue "D41" "0.007354"
Set-TextValue "E41" "-1.36%"

Set-TextValue "D42" "0.1305"
Set-TextValue "E42" "-4.04%"

Set-TextValue "D43" "0.002010"
Set-TextValue "E43" "-15.96%"

Set-TextValue "E44" "-0.84%"

Set-TextValue "D45" "0.00005837"
Set-TextValue "E45" "-2.46%"

Set-TextValue "D46" "0.00000000752"
Set-TextValue "E46" "0.19%"

Set-TextValue "D47" "2.200"
Set-TextValue "E47" "168.16%"

Set-TextValue "D49" "0.00002105"
Set-TextValue "E49" "0.19%"

Set-TextValue "D50" "0.0002005"
Set-TextValue "E50" "0.19%"
